$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.196.04'
$ws.Range('E2').Value = '  +3.14%  '
$ws.Range('D3').Value = '2.633.32'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.10'
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.00'
$ws.Range('E6').Value = '  +4.15%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('E9').Value = '  +8.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.399'
$ws.Range('E10').Value = '  +4.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.80'
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('E12').Value = '  +1.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.97'
$ws.Range('E13').Value = '  +5.88%  '
$ws.Range('E14').Value = '  +18.78%  '
$ws.Range('D15').Value = '3.103.35'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '65.078.79'
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').Value = '2.617.00'
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('E18').Value = '  +1.98%  '
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '354.04'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('E21').Value = '  +5.60%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.93'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.71'
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.47'
$ws.Range('E25').Value = '  +2.54%  '
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.12'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.163'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('D29').Value = '0.0₃0950'
$ws.Range('E29').Value = '  +12.99%  '
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '522.56'
$ws.Range('E31').Value = '  -6.52%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.10'
$ws.Range('E32').Value = '  +4.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.79'
$ws.Range('E33').Value = '  +2.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.66'
$ws.Range('E34').Value = '  +8.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.33'
$ws.Range('E35').Value = '  +3.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.426'
$ws.Range('E36').Value = '  +4.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.04'
$ws.Range('E37').Value = '  +6.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '164.75'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '20.19'
$ws.Range('E39').Value = '  +3.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.24'
$ws.Range('E42').Value = '  +6.78%  '
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('E44').Value = '  +3.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0614'
$ws.Range('E45').Value = '  +5.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.96'
$ws.Range('E46').Value = '  +1.72%  '
$ws.Range('E47').Value = '  +9.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.647'
$ws.Range('E48').Value = '  +2.99%  '
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0985'
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.41'
$ws.Range('E51').Value = '  +2.47%  '
